# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values were recalculated for data rows 2-7.
$kValues = @{
    2 = 0
    3 = 1
    4 = 0
    5 = 1
    6 = 0
    7 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
